$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2 ("Лист2"): new Timer5/Timer6 + ESila/VE8/VK014 base-address table
# Values are written in the exact order new distinct strings first appear so
# that the shared-string table is built up with the same indices as the
# target workbook.
# ---------------------------------------------------------------------------

# Row 9 - table header for the 3 new columns (BE8 / VK014 / ESila)
$ws2.Range("C9").Value = "BE8"
$ws2.Range("D9").Value = "ВК014"
$ws2.Range("E9").Value = "Esila"

# Column E (ESila) base addresses for Timer1..Timer3
$ws2.Range("E10").Value = "0x4009_4000"
$ws2.Range("E11").Value = "0x4009_5000"
$ws2.Range("E12").Value = "0x4009_6000"

# Column C (BE8, mirrored into D for VK014) base addresses for Timer1..Timer4
$ws2.Range("C10").Value = "0x4008_A000"
$ws2.Range("C11").Value = "0x4008_B000"
$ws2.Range("C12").Value = "0x4008_C000"
$ws2.Range("C13").Value = "0x4008_D000"

# Column D (VK014) base addresses for the extra Timer5 / Timer6 rows
$ws2.Range("D14").Value = "0x4008_E000"
$ws2.Range("D15").Value = "0x4008_F000"

# New timer row labels
$ws2.Range("B14").Value = "Timer5"
$ws2.Range("B15").Value = "Timer6"

# Footnote caption (merged B16:E16)
$ws2.Range("B16").Value = "Timers аналогичны 1986ВЕ1Т"

# --- remaining cells that reuse already-existing shared strings -----------
$ws2.Range("B10").Value = "Timer1"
$ws2.Range("B11").Value = "Timer2"
$ws2.Range("B12").Value = "Timer3"
$ws2.Range("B13").Value = "Timer4"

$ws2.Range("D10").Value = "0x4008_A000"
$ws2.Range("D11").Value = "0x4008_B000"
$ws2.Range("D12").Value = "0x4008_C000"
$ws2.Range("D13").Value = "0x4008_D000"

$ws2.Range("E13").Value = "-"
$ws2.Range("C14").Value = "-"
$ws2.Range("E14").Value = "-"
$ws2.Range("C15").Value = "-"
$ws2.Range("E15").Value = "-"

Write-Output "values done"

# ---------------------------------------------------------------------------
# Styles
# ---------------------------------------------------------------------------

# s=15 (default font, thin box border) - reuse the existing style used by the
# rest of the table (e.g. B3) via format-only paste so no new cellXf is
# minted.
$ws2.Range("B3").Copy()
$ws2.Range("B9").PasteSpecial(-4122)
$ws2.Range("B10:E12").PasteSpecial(-4122)
$ws2.Range("B13:D13").PasteSpecial(-4122)
$ws2.Range("B14").PasteSpecial(-4122)
$ws2.Range("D14").PasteSpecial(-4122)
$ws2.Range("B15").PasteSpecial(-4122)
$ws2.Range("D15").PasteSpecial(-4122)

# s=26 (bold font, thin box border) - header-style cells, reuse C1's style.
$ws2.Range("C1").Copy()
$ws2.Range("C9:E9").PasteSpecial(-4122)

# s=27 (default font, thin box border, quote-prefixed "-" centered) - reuse
# G5's style.
$ws2.Range("G5").Copy()
$ws2.Range("E13").PasteSpecial(-4122)
$ws2.Range("C14").PasteSpecial(-4122)
$ws2.Range("E14").PasteSpecial(-4122)
$ws2.Range("C15").PasteSpecial(-4122)
$ws2.Range("E15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Output "styles (existing) done"

# s=28 (new bold+italic font, top-only thin border, centered) - build once on
# a single cell so only one new font/border/cellXf is minted, then propagate
# via format-only paste, then merge.
$capCell = $ws2.Range("B16")
$capCell.Font.Bold = $true
$capCell.Font.Italic = $true
$capCell.HorizontalAlignment = -4108
$capCell.Borders.Item(8).LineStyle = 1

$ws2.Range("B16").Copy()
$ws2.Range("C16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B16:E16").Merge()

Write-Output "styles (new) done"

# ---------------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------------
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 108
$win.ScrollColumn = 1

$ws2.Activate()
$ws2.Range("E21").Select()

Write-Output "done"
